# TOM 14 "Ich bin in einer Tomate" - add cross-references to other training
# cards (TOM 11/12 and SCR 10-14) in the "Inhaltsplatzhalter 3" body text box
# of slide 1.

$p  = $ppt.ActivePresentation
$sl = $p.Slides.Item(1)
$sh = $sl.Shapes.Item(3)          # "Inhaltsplatzhalter 3"
$tr = $sh.TextFrame.TextRange

# ---------------------------------------------------------------------
# Change 1: "...der Pomodoro    nicht unmittelbar verfolgen.<nbsp>"
#           four leading spaces -> a single leading space.
# The run ends with a non-breaking space (U+00A0); it must stay untouched
# (it gets corrupted if round-tripped through a PowerShell string read),
# so we rebuild it from a literal [char]0x00A0 instead of reusing the
# value read back from .Text.
# ---------------------------------------------------------------------
$full = $tr.Text
$needle = "nicht unmittelbar verfolgen"
$idx0 = $full.IndexOf($needle)            # 0-based index of "nicht..."
$runStart0 = $idx0 - 4                    # run starts 4 spaces earlier
$vIdx0 = $full.IndexOf("verfolgen.")
$runEnd0 = $vIdx0 + "verfolgen.".Length + 1   # include trailing nbsp char
$runLen = $runEnd0 - $runStart0

$run1 = $tr.Characters($runStart0 + 1, $runLen)
$nbsp = [char]0x00A0
$run1.Text = " nicht unmittelbar verfolgen." + $nbsp

# ---------------------------------------------------------------------
# Change 2: add references to other training cards in the paragraph
# that talks about the Tagesplan / Backlog.
# ---------------------------------------------------------------------

# 2a) Shrink the run "werden muss, ... willst, kann in Dein " down to
#     "werden muss, ... Tagesplan (TOM 11, TOM 12). Was du später erledigen "
#     and drop its "willst, kann in Dein " tail (moved to a new run below).
$full = $tr.Text
$runBText = "werden muss, kommt sofort in den Tagesplan. Was du später erledigen willst, kann in Dein "
$idxB0 = $full.IndexOf($runBText)
$runB = $tr.Characters($idxB0 + 1, $runBText.Length)
$runB.Text = "werden muss, kommt sofort in den Tagesplan (TOM 11, TOM 12). Was du später erledigen "

# 2b) Insert a new run "willst, kann in Dein " right before "Backlog"
#     (after the existing line break), then re-assert "Backlog" itself so
#     it becomes its own run again (matching the original run boundary).
$full = $tr.Text
$backlogIdx0 = $full.IndexOf("Backlog")
$backlogRange = $tr.Characters($backlogIdx0 + 1, 7)
$backlogRange.InsertBefore("willst, kann in Dein ")

$full = $tr.Text
$backlogIdx0b = $full.IndexOf("Backlog")
$backlogRange2 = $tr.Characters($backlogIdx0b + 1, 7)
$backlogRange2.Text = "Backlog"

# 2c) Extend the trailing "." run after "Backlog" into the full
#     reference sentence.
$full = $tr.Text
$backlogIdx0c = $full.IndexOf("Backlog")
$dotPos0 = $backlogIdx0c + 7
$dotRange = $tr.Characters($dotPos0 + 1, 1)
$dotRange.Text = " (SCR 10-12) oder Deinen Wochenplan (SCR 13, SCR 14)."
